# Update the "想去人数" (want-to-go count) figures in column F of the
# "展览" and "全部类型" worksheets to reflect the newly scraped values.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# Sheet "展览": row -> new value
$exhibitionUpdates = @{
    3  = 566
    4  = 211
    6  = 518
    7  = 112
    8  = 127
    9  = 51
    10 = 6934
    12 = 388
    13 = 3324
    14 = 228
    15 = 416
    17 = 571
    18 = 44
}

foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型": row -> new value
$allTypesUpdates = @{
    5  = 566
    6  = 211
    8  = 518
    9  = 112
    10 = 127
    11 = 51
    13 = 6934
    16 = 388
    17 = 3324
    18 = 228
    19 = 416
    21 = 571
    22 = 44
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}

$wb.Save()
